# ISYS3001.docx - "Making change for the assignment"
#
# 1. The placeholder "…" paragraph becomes a real comment about the edit.
# 2. The trailing empty paragraph gets the author's GitHub user id, typed
#    as two separate runs (the label, then the id) just like the rest of
#    the document's hand-edited paragraphs.

$d = $word.ActiveDocument

# 1) Replace the ellipsis placeholder with the real commit comment.
$d.Content.Find.Execute("…", $true, $false, $false, $false, $false, `
                         $true, 1, $false, `
                         "Changing this file for the assignment", 2) | Out-Null

# 2) Turn the final (empty) paragraph into the GitHub user id line, typed as
#    two distinct runs. Briefly turning on TrackRevisions forces Word to
#    keep each insertion as its own run instead of silently coalescing them
#    with whatever is already in the paragraph; AcceptAllRevisions then
#    folds the tracked insertions back into normal (untracked) text.
$d.TrackRevisions = $true

$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertAfter("My github user id: ")

$idRange = $d.Paragraphs.Last.Range
$idRange.Collapse(0)
$idRange.InsertAfter("karan90341")

$d.TrackRevisions = $false
$d.AcceptAllRevisions() | Out-Null
